$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Label" header in H1, reusing the same header style as the other
# header cells (bold / bordered / centered) by copying formats from G1.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Updated D/E/F values (refit results) and new H (Label) column values
$updates = @(
    @{ Row = 2;  D = 0.5402005174746579; E = 0.5402005174746579; H = 0 },
    @{ Row = 3;  D = 0.6012545495941969; E = 0.6012545495941969; H = 0 },
    @{ Row = 4;  D = 0.4219870694208052; E = 0.4219870694208052; H = 0 },
    @{ Row = 5;  D = 0.407899665454832;  E = 0.407899665454832;  H = 0 },
    @{ Row = 6;  D = 0.4948077986537394; E = 0.4948077986537394; H = 0 },
    @{ Row = 7;  D = 0.4138238637374368; E = 0.5861761362625632; H = 1 },
    @{ Row = 8;  D = 0.5468922917754426; E = 0.4531077082245574; H = 1 },
    @{ Row = 9;  D = 0.9864197003452576; E = 0.01358029965474239; H = 1 },
    @{ Row = 10; D = 0.5880380439698119; E = 0.4119619560301881; H = 1 },
    @{ Row = 11; D = 0.9585630925011284; E = 0.04143690749887163; F = 0.5524226427078247; H = 1 },
    @{ Row = 12; H = 0 },
    @{ Row = 13; H = 0 },
    @{ Row = 14; H = 0 },
    @{ Row = 15; H = 0 },
    @{ Row = 16; H = 0 },
    @{ Row = 17; H = 1 },
    @{ Row = 18; H = 1 },
    @{ Row = 19; H = 1 },
    @{ Row = 20; H = 1 },
    @{ Row = 21; H = 1 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $u.F }
    if ($u.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $u.H }
}
